# Update metricas_recorrencia_anual row 8 (year 2025) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1200
$ws.Range("D8").Value = 195
$ws.Range("E8").Value = 1005
$ws.Range("F8").Value = 7.998359310910582
$ws.Range("G8").Value = 83.75
$ws.Range("H8").Value = 16.25
